# Generate Report for Handoff
#
# A new handoff pass ran for the four "Ready for handoff" rows
# (1d1d2628-…, 65ee4d5e-…, 7ca21cb1-…, f6614653-…) on both the zh-cn and
# de-de localization-status sheets:
#   - Priority moved from "low" to "ht" for those rows.
#   - The "Latest Handoff Datetime" for those rows was refreshed to the
#     new handoff-generation timestamp (per-sheet).

$wb = $excel.ActiveWorkbook

$rows = 4, 5, 6, 7

$sheetInfo = @{
    "zh-cn" = "2016-08-16 18:26:57"
    "de-de" = "2016-08-16 18:27:06"
}

foreach ($sheetName in $sheetInfo.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $handoffTime = $sheetInfo[$sheetName]

    foreach ($r in $rows) {
        $ws.Range("E$r").Value = "ht"
        $ws.Range("H$r").Value = $handoffTime
    }
}

# The Overview sheet's "Latest HO Xliff Generate Date" column (G) tracks the
# most recent handoff timestamp across all languages for each file. With the
# de-de handoff now the latest (2016-08-16 18:27:06), rows 4-7 pick it up too.
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G$r").Value = $sheetInfo["de-de"]
}
